# Apply the "Add files via upload" re-upload edit:
#   - rename the two worksheet tabs
#   - make the second sheet ("steel") the active/selected tab instead of the first
#
# (The source diff also shows styles.xml's custom number format being
#  renumbered from numFmtId 166 to 164 with no change to its formatCode
#  ("0.0") or to any cell's style index -- a pure internal re-ID that
#  Excel's save pass performs on its own and that isn't reachable through
#  the Number/Styles object model, so there is nothing for this script to
#  call for that part; the format itself keeps working unchanged.)

$wb = $excel.ActiveWorkbook

$wsAluminium = $wb.Worksheets.Item(1)
$wsSteel     = $wb.Worksheets.Item(2)

$wsAluminium.Name = "aluminium"
$wsSteel.Name = "steel"

# Switch the active tab from the first sheet to the second one.
$wsSteel.Activate()
